$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-8 from 2023-09-01 (45170)
# to 2023-09-05 (45174), preserving existing cell formatting.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
